$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "test"
$ws.Range("C2").Value = [double]"2.71403116699935"
$ws.Range("D2").Value = [double]"53.5997396975144"
$ws.Range("E2").Value = [double]"31.85301908166584"
$ws.Range("F2").Value = [double]"4.044213102507071"
$ws.Range("G2").Value = [double]"0.817380984054707"
$ws.Range("H2").Value = [double]"0.04904287205040092"
$ws.Range("I2").Value = [double]"6.885383400118287"
$ws.Range("J2").Value = [double]"0.006959032432913243"
$ws.Range("K2").Value = [double]"0.002917806098027364"
$ws.Range("M2").Value = [double]"9.69714560178551e-05"
$ws.Range("N2").Value = [double]"0.02720710064251893"
$ws.Range("O2").Value = [double]"8.681164389933109e-06"
$ws.Range("P2").Value = [double]"4.100164157773595e-09"
$ws.Range("Q2").Value = [double]"1.679633603058436e-09"
$ws.Range("R2").Value = [double]"1.011498257588196e-08"
$ws.Range("S2").Value = [double]"3.314433459902785e-08"
$ws.Range("T2").Value = [double]"5.236005979697833e-08"
$ws.Range("U2").Value = [double]"1.839843026964226e-09"
$ws.Range("V2").Value = [double]"5.387692600459592e-11"
$ws.Range("W2").Value = [double]"7.693327657563124e-14"
$ws.Range("X2").Value = [double]"3.733564744705091e-14"
$ws.Range("Y2").Value = [double]"3.982510026496743e-15"
$ws.Range("Z2").Value = [double]"1.990720783505968e-12"
$ws.Range("AA2").Value = [double]"6.4352727547928e-13"
$ws.Range("AB2").Value = [double]"1.088319412220221e-13"
$ws.Range("AC2").Value = [double]"2.904570109423545e-13"
$ws.Range("AD2").Value = [double]"3.086092124408925e-14"
$ws.Range("AF2").Value = [double]"1.926528414258276e-16"
$ws.Range("AG2").Value = [double]"1.379464203440023e-17"
$ws.Range("AH2").Value = [double]"2.077473365773284e-20"
$ws.Range("AI2").Value = [double]"5.685688199576031e-19"
$ws.Range("AJ2").Value = [double]"1.206445672521582e-17"
$ws.Range("AL2").Value = [double]"6.279493501439638e-21"
$ws.Range("AM2").Value = [double]"1.299115057040853e-22"
$ws.Range("AN2").Value = [double]"3.16846952943607e-27"
$ws.Range("AR2").Value = [double]"0.02010164240889484"
$ws.Range("AS2").Value = [double]"2.199794416385694"
$ws.Range("AT2").Value = [double]"4.453779597059164"
$ws.Range("AU2").Value = [double]"2.006625339165485"
$ws.Range("AV2").Value = [double]"1.027804084269509"
$ws.Range("AW2").Value = [double]"0.1976994101726686"
$ws.Range("AX2").Value = [double]"88.16960308249729"
$ws.Range("AY2").Value = [double]"0.05562716627717863"
$ws.Range("AZ2").Value = [double]"0.05204657225355335"
$ws.Range("BB2").Value = [double]"0.0044696754182013"
$ws.Range("BC2").Value = [double]"1.812051816575728"
$ws.Range("BD2").Value = [double]"0.0002757755083410862"
$ws.Range("BE2").Value = [double]"1.585287238625901e-06"
$ws.Range("BF2").Value = [double]"7.49972886482383e-07"
$ws.Range("BG2").Value = [double]"1.038411247960721e-05"
$ws.Range("BH2").Value = [double]"4.053768177584842e-05"
$ws.Range("BI2").Value = [double]"6.484807992945085e-05"
$ws.Range("BJ2").Value = [double]"3.258966373266431e-06"
$ws.Range("BK2").Value = [double]"4.914817523177494e-08"
$ws.Range("BL2").Value = [double]"1.741136525713667e-10"
$ws.Range("BM2").Value = [double]"4.341072097176452e-11"
$ws.Range("BN2").Value = [double]"1.356188171798649e-11"
$ws.Range("BO2").Value = [double]"5.293435375984415e-09"
$ws.Range("BP2").Value = [double]"1.727677512580947e-09"
$ws.Range("BQ2").Value = [double]"6.072815404192162e-10"
$ws.Range("BR2").Value = [double]"6.923832348243797e-10"
$ws.Range("BS2").Value = [double]"2.054013591516476e-10"
$ws.Range("BU2").Value = [double]"1.914755879560026e-12"
$ws.Range("BV2").Value = [double]"9.617625506096119e-14"
$ws.Range("BW2").Value = [double]"2.387024499796258e-16"
$ws.Range("BX2").Value = [double]"5.107985078971899e-15"
$ws.Range("BY2").Value = [double]"1.505166662102754e-13"
$ws.Range("CA2").Value = [double]"2.141165158759857e-16"
$ws.Range("CB2").Value = [double]"5.210042853027194e-18"
$ws.Range("CC2").Value = [double]"3.779032516856507e-22"
$ws.Range("CG2").Value = [double]"0.0003343561715082666"
$ws.Range("CH2").Value = [double]"0.0216545686623432"
$ws.Range("CI2").Value = [double]"0.04106417206088502"
$ws.Range("CJ2").Value = [double]"0.02686938134111512"
$ws.Range("CK2").Value = [double]"0.02516489392536681"
$ws.Range("CL2").Value = [double]"0.01839440068045068"
$ws.Range("CM2").Value = [double]"69.52393312163424"
$ws.Range("CN2").Value = [double]"0.02190942470577657"
$ws.Range("CO2").Value = [double]"0.05334476560554074"
$ws.Range("CQ2").Value = [double]"0.02340469005358193"
$ws.Range("CR2").Value = [double]"30.23598174273442"
$ws.Range("CS2").Value = [double]"0.002289249307558654"
$ws.Range("CT2").Value = [double]"2.627030976199887e-05"
$ws.Range("CU2").Value = [double]"1.444835178373056e-05"
$ws.Range("CV2").Value = [double]"0.0003766982814818166"
$ws.Range("CW2").Value = [double]"0.001809976787624915"
$ws.Range("CX2").Value = [double]"0.003131012814064291"
$ws.Range("CY2").Value = [double]"0.0002910333401375669"
$ws.Range("CZ2").Value = [double]"2.326505472742829e-06"
$ws.Range("DA2").Value = [double]"3.507543645587063e-08"
$ws.Range("DB2").Value = [double]"8.942368410339862e-09"
$ws.Range("DC2").Value = [double]"5.844827793614924e-09"
$ws.Range("DD2").Value = [double]"1.69278427433024e-06"
$ws.Range("DE2").Value = [double]"6.32286678097174e-07"
$ws.Range("DF2").Value = [double]"5.079505069542724e-07"
$ws.Range("DG2").Value = [double]"2.734114140146589e-07"
$ws.Range("DH2").Value = [double]"3.023411311194524e-07"
$ws.Range("DJ2").Value = [double]"6.661682463367941e-09"
$ws.Range("DK2").Value = [double]"2.850277848176648e-10"
$ws.Range("DL2").Value = [double]"1.490582537584066e-12"
$ws.Range("DM2").Value = [double]"3.397163254437789e-11"
$ws.Range("DN2").Value = [double]"1.096063059062184e-09"
$ws.Range("DP2").Value = [double]"1.145005110717641e-11"
$ws.Range("DQ2").Value = [double]"5.512801110892305e-13"
$ws.Range("DR2").Value = [double]"5.927978082991556e-16"
$ws.Range("DY2").Value = [double]"9.548813714624738e-24"
$ws.Range("DZ2").Value = [double]"2.634031516911138e-16"
$ws.Range("EA2").Value = [double]"1.958975334249318e-09"
$ws.Range("EB2").Value = [double]"0.2170734125160322"
$ws.Range("EC2").Value = [double]"0.0001852819406223197"
$ws.Range("ED2").Value = [double]"0.003156490882481011"
$ws.Range("EF2").Value = [double]"0.008358894842469388"
$ws.Range("EG2").Value = [double]"43.82403759677042"
$ws.Range("EH2").Value = [double]"0.008094919327571622"
$ws.Range("EI2").Value = [double]"0.01001860694413209"
$ws.Range("EJ2").Value = [double]"0.01002073680527772"
$ws.Range("EK2").Value = [double]"1.187065851454961"
$ws.Range("EL2").Value = [double]"9.270228871785305"
$ws.Range("EM2").Value = [double]"20.25061000358655"
$ws.Range("EN2").Value = [double]"8.295488312561353"
$ws.Range("EO2").Value = [double]"0.01454670895244267"
$ws.Range("EP2").Value = [double]"0.006999662110418134"
$ws.Range("EQ2").Value = [double]"0.00274975897485704"
$ws.Range("ER2").Value = [double]"0.008166278399896946"
$ws.Range("ES2").Value = [double]"1.144445582951756"
$ws.Range("ET2").Value = [double]"0.5996381060692174"
$ws.Range("EU2").Value = [double]"3.03785603396554"
$ws.Range("EV2").Value = [double]"0.337150637334813"
$ws.Range("EW2").Value = [double]"7.759131478923479"
$ws.Range("EY2").Value = [double]"1.238941282489459"
$ws.Range("EZ2").Value = [double]"0.04083139569952442"
$ws.Range("FA2").Value = [double]"0.001166611370331441"
$ws.Range("FB2").Value = [double]"0.02683206129718998"
$ws.Range("FC2").Value = [double]"1.048783509875114"
$ws.Range("FE2").Value = [double]"1.206276225858535"
$ws.Range("FF2").Value = [double]"0.3091519853225819"
$ws.Range("FG2").Value = [double]"0.1329936990286836"
